# Reproduce the "Add files via upload" edit to ML/fEX_scorebook.xlsx.
#
# The refreshed data replaces the repeated "0.671313 / 0.652727 / 0.6527272727272727"
# block that filled rows 20-25 with new, row-specific Train/Test/Precision/Recall
# text values, updates the TP/TN/FP/FN counters that go with them, bumps B20:B25
# from 0.3333 to 0.5, refreshes row 26's Precision/Recall text, and moves the
# sheet's active selection to B25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even though it looks like a
# number (matches the source file, where these columns hold shared-string text
# such as "0.579403" rather than numeric cells). Forcing the cell to Text format
# before the write keeps Excel from re-interpreting the literal as a number;
# switching back to the Normal style afterwards drops the explicit format again
# so the cell itself carries no leftover style index.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---- Row 20 ----
$ws.Range("B20").Value = 0.5
Set-TextValue "E20" "0.579403"
Set-TextValue "F20" "0.622642"
$ws.Range("G20").Value = 263
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 160
Set-TextValue "K20" "0.6226415094339622"
Set-TextValue "L20" "0.6226415094339622"

# ---- Row 21 ----
$ws.Range("B21").Value = 0.5
Set-TextValue "E21" "0.504667"
Set-TextValue "F21" "0.486000"
$ws.Range("H21").Value = 242
$ws.Range("I21").Value = 257
Set-TextValue "K21" "0.4866415094339622"
Set-TextValue "L21" "0.4866415094339622"

# ---- Row 22 ----
$ws.Range("B22").Value = 0.5
Set-TextValue "E22" "0.504667"
Set-TextValue "F22" "0.486000"
$ws.Range("H22").Value = 242
$ws.Range("I22").Value = 257
Set-TextValue "K22" "0.4866415094339622"
Set-TextValue "L22" "0.4866415094339622"

# ---- Row 23 ----
$ws.Range("B23").Value = 0.5
Set-TextValue "E23" "0.504667"
Set-TextValue "F23" "0.486000"
$ws.Range("H23").Value = 242
$ws.Range("I23").Value = 257
Set-TextValue "K23" "0.4866415094339622"
Set-TextValue "L23" "0.4866415094339622"

# ---- Row 24 ----
$ws.Range("B24").Value = 0.5
Set-TextValue "E24" "0.504667"
Set-TextValue "F24" "0.486000"
$ws.Range("H24").Value = 242
$ws.Range("I24").Value = 257
Set-TextValue "K24" "0.4866415094339622"
Set-TextValue "L24" "0.4866415094339622"

# ---- Row 25 ----
$ws.Range("B25").Value = 0.5
Set-TextValue "E25" "0.506346"
Set-TextValue "F25" "0.478000"
$ws.Range("H25").Value = 260
$ws.Range("I25").Value = 239
Set-TextValue "K25" "0.6527272727272727"
Set-TextValue "L25" "0.6527272727272727"

# ---- Row 26 (only Precision/Recall text refreshed) ----
Set-TextValue "K26" "0.4787272727272727"
Set-TextValue "L26" "0.4787272727272727"

# ---- Selection moves to B25 ----
$ws.Range("B25").Select()
